$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the percentage values (row 4-6, columns C:E)
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Value = 32.479999999999997
$ws.Range("E4").Value = 38.92

$ws.Range("C5").Value = 41.17
$ws.Range("D5").Value = 0.03
$ws.Range("E5").Value = 57.56

$ws.Range("C6").Value = 33.94
$ws.Range("D6").Value = 67.52
$ws.Range("E6").Value = 0.04

# Update the active selection cell to match the author's final cursor position
$ws.Range("G10").Select()
